# "Generate Report for Handoff"
#
# The file 68d3c360-dbcf-4407-aab7-bb8a84661ca2.md has just been handed
# off again, so its status moves from "Handed back: in sync with en-US"
# to "Ready for handoff" everywhere it's reported, and the per-locale
# "Latest Handoff Datetime" is stamped with the new handoff time.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is 68d3c360-dbcf-4407-aab7-bb8a84661ca2.md;
# its zh-cn and de-de status columns both flip to "Ready for handoff".
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the same file's zh-cn handoff record.
# Status -> "Ready for handoff"; Latest Handoff Datetime -> new timestamp.
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "2016-03-10 04:57:34"

# de-de sheet: row 3 is the same file's de-de handoff record.
# Status -> "Ready for handoff"; Latest Handoff Datetime -> new timestamp.
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "2016-03-10 04:57:38"
